# Natmi following Dr Hou advice
# Rewrites the LR-pairs (Nts-Sort1) results sheet:
#   - recomputes the Ligand-expressing-cells / specificity / weight metrics
#     for the existing ECs -> {ECs,FAPs,M1,M2,sCs} rows (rows 2-6)
#   - adds the new sCs -> {ECs,FAPs,M1,M2,sCs} sending-cluster rows (rows 7-11)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
  @("ECs", "Nts", "Sort1", "ECs",  3, 1, 1.747141666666667, 5.241425,   0.9598684458717283,  0.9598684458717284,  3, 1, 0.72097,            2.16291,            0.03318293102712134, 0.03318293102712135, 1.259636727416667, 11.33673054675,    0.03185124843447171, 0.03185124843447172),
  @("ECs", "Nts", "Sort1", "FAPs", 3, 1, 1.747141666666667, 5.241425,   0.9598684458717283,  0.9598684458717284,  3, 1, 1.851391,           5.554173,           0.08521100719479757, 0.08521100719479759, 3.234642357391666, 29.111781216525,   0.081791357047235,   0.08179135704723503),
  @("ECs", "Nts", "Sort1", "M1",   3, 1, 1.747141666666667, 5.241425,   0.9598684458717283,  0.9598684458717284,  3, 1, 8.824443666666665, 26.473331,          0.4061485298191571,  0.4061485298191572,  15.41755321518611, 138.757978936675,  0.3898491581106016,  0.3898491581106018),
  @("ECs", "Nts", "Sort1", "M2",   3, 1, 1.747141666666667, 5.241425,   0.9598684458717283,  0.9598684458717284,  3, 1, 8.620174666666665, 25.860524,          0.3967469678429597,  0.3967469678429598,  15.06066633407777, 135.5459970067,    0.3808248954277423,  0.3808248954277424),
  @("ECs", "Nts", "Sort1", "sCs",  3, 1, 1.747141666666667, 5.241425,   0.9598684458717283,  0.9598684458717284,  3, 1, 1.710155,           5.130465,           0.07871056411596419, 0.07871056411596419, 2.987883056958333, 26.890947512625,   0.07555178685167757, 0.07555178685167759),
  @("sCs", "Nts", "Sort1", "ECs",  1, 0.3333333333333333, 0.073047, 0.219141, 0.04013155412827168, 0.04013155412827169, 3, 1, 0.72097,            2.16291,            0.03318293102712134, 0.03318293102712135, 0.05266469559,      0.47398226031,       0.001331682592649626, 0.001331682592649626),
  @("sCs", "Nts", "Sort1", "FAPs", 1, 0.3333333333333333, 0.073047, 0.219141, 0.04013155412827168, 0.04013155412827169, 3, 1, 1.851391,           5.554173,           0.08521100719479757, 0.08521100719479759, 0.135238558377,     1.217147025393,      0.003419650147562567, 0.003419650147562568),
  @("sCs", "Nts", "Sort1", "M1",   1, 0.3333333333333333, 0.073047, 0.219141, 0.04013155412827168, 0.04013155412827169, 3, 1, 8.824443666666665, 26.473331,          0.4061485298191571,  0.4061485298191572,  0.6445991365189999, 5.801392228670999,   0.01629937170855547,  0.01629937170855548),
  @("sCs", "Nts", "Sort1", "M2",   1, 0.3333333333333333, 0.073047, 0.219141, 0.04013155412827168, 0.04013155412827169, 3, 1, 8.620174666666665, 25.860524,          0.3967469678429597,  0.3967469678429598,  0.6296778988759999, 5.667101089884,      0.0159220724152174,   0.01592207241521741),
  @("sCs", "Nts", "Sort1", "sCs",  1, 0.3333333333333333, 0.073047, 0.219141, 0.04013155412827168, 0.04013155412827169, 3, 1, 1.710155,           5.130465,           0.07871056411596419, 0.07871056411596419, 0.124921692285,     1.124295230565,      0.003158777264286616, 0.003158777264286617)
)

for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $rows[$i]
    for ($j = 0; $j -lt $r.Length; $j++) {
        $ws.Cells.Item($i + 2, $j + 1).Value = $r[$j]
    }
}
